$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used for each data row (B,C,D,E,F skip G,H then I,J,K,L,M,N)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Row data: each inner array holds values for columns in $cols, for sheet rows 2..25
$rowData = @{
    2 = @(1.02, 1.093577141556132, 1.085909644096236, 1.104818817500664, 1.10642821456568, 1.050679878400681, 1.098392401064319, 1.088568854182986, 1.107430037964381, 1.109035432101395, 1.0999522442612)
    3 = @(1.02, 1.095540298482043, 1.087424598182881, 1.106713533060759, 1.108256164959711, 1.051138500537724, 1.100017202219391, 1.089901417140266, 1.109144990456562, 1.110684067663543, 1.101579352820279)
    4 = @(1.02, 1.096807235511342, 1.088401641174331, 1.107936503054189, 1.109435906414277, 1.051432454490827, 1.101064899575203, 1.09075993146114, 1.110251198271932, 1.111747323285974, 1.102628538026509)
    5 = @(1.02, 1.097339071755595, 1.088811630180363, 1.10844993035217, 1.109931153218401, 1.051555366476053, 1.10150449213988, 1.091119968270954, 1.11071543144351, 1.112193488944654, 1.10306875486305)
    6 = @(1.02, 1.097428323986287, 1.088880424991255, 1.108536095979647, 1.110014265816147, 1.05157596502316, 1.101578251738847, 1.091180368667799, 1.110793330835121, 1.112268354091591, 1.103142619209107)
    7 = @(1.02, 1.096814344987596, 1.08840712243378, 1.107943366257446, 1.109442526718249, 1.051434099457425, 1.101070776787301, 1.09076474574028, 1.11025740456406, 1.111753288200917, 1.102634423584922)
    8 = @(1.02, 1.094241306486111, 1.086422308089348, 1.105459785311805, 1.10704662391651, 1.05083545596259, 1.09894227738326, 1.089019983240834, 1.108010345134139, 1.109593335179741, 1.100502901467586)
    9 = @(1.02, 1.089680693811884, 1.082899402208992, 1.101059324212715, 1.102800489764309, 1.049758844325521, 1.095162843395329, 1.085916213428692, 1.104023320303172, 1.105759521215947, 1.096718100249884)
    10 = @(1.02, 1.086621169558001, 1.080532808715057, 1.098108359976055, 1.099952340982138, 1.049026157356718, 1.092622863038337, 1.083826494715287, 1.101345803456547, 1.103183997532494, 1.094174512829446)
    11 = @(1.02, 1.085291550571647, 1.079503577622787, 1.096826196217795, 1.098714691568898, 1.048705274629445, 1.091517964035012, 1.082916560006265, 1.100181545215609, 1.102063877862033, 1.093068044742802)
    12 = @(1.02, 1.084796922006525, 1.079120584601654, 1.096349263878866, 1.098254293929797, 1.048585533317561, 1.091106773507705, 1.082577790962741, 1.099748334646454, 1.101647059352054, 1.092656270277782)
    13 = @(1.02, 1.084903055789405, 1.079202769448805, 1.096451598596587, 1.098353081879553, 1.048611243285163, 1.091195010975736, 1.082650493583632, 1.099841294203851, 1.10173650287852, 1.092744633053139)
    14 = @(1.02, 1.085250679836864, 1.079471933494394, 1.096786786834828, 1.098676648942264, 1.048695388057935, 1.091483990963825, 1.08288857323543, 1.100145751380918, 1.10202943905598, 1.093034023425957)
    15 = @(1.02, 1.085464762529425, 1.07963768228193, 1.096993216648242, 1.09887591860509, 1.04874715919246, 1.091661936867077, 1.083035158268806, 1.100333237063616, 1.102209825901155, 1.093212222032802)
    16 = @(1.02, 1.086709308296518, 1.080601019495058, 1.09819335887102, 1.100034385569339, 1.049047376358315, 1.092696083008541, 1.083886775791759, 1.101422966954641, 1.103258231237904, 1.094247836780405)
    17 = @(1.02, 1.08748867039216, 1.081204083552457, 1.098944989101129, 1.100759873193946, 1.049234719624274, 1.093343403925394, 1.08441960316075, 1.102105205817583, 1.10391454213453, 1.094896076967209)
    18 = @(1.02, 1.087942795230228, 1.081555409513492, 1.099382981507624, 1.101182616475863, 1.049343644530046, 1.093720487375789, 1.084729904042959, 1.102502675172604, 1.104296885664707, 1.095273695919356)
    19 = @(1.02, 1.088097562102908, 1.081675129994775, 1.099532254889807, 1.101326690174152, 1.049380726068707, 1.093848980935198, 1.084835626300305, 1.102638123017376, 1.104427175549053, 1.095402371954358)
    20 = @(1.02, 1.087405100370104, 1.081139425204276, 1.098864389974226, 1.100682078975199, 1.049214655637521, 1.093274003126143, 1.084362486454728, 1.102032056707498, 1.103844175100638, 1.094826577610863)
    21 = @(1.02, 1.08514833408599, 1.079392690602919, 1.096688101192788, 1.098581385436957, 1.048670624788703, 1.091398915304674, 1.082818486344892, 1.100056117290688, 1.101943197688477, 1.092948826949612)
    22 = @(1.02, 1.083725070697813, 1.078290443770972, 1.095315838788274, 1.097256653739682, 1.048325378687608, 1.090215439570485, 1.081843198017558, 1.098809397415801, 1.100743592617261, 1.091763670544134)
    23 = @(1.02, 1.084479989483212, 1.078875151286252, 1.096043682270997, 1.097959299897495, 1.048508705015515, 1.090843258932372, 1.082360650322769, 1.099470728254361, 1.101379948326191, 1.092392381481526)
    24 = @(1.02, 1.087442863498668, 1.081168642860376, 1.098900810550079, 1.100717232135632, 1.049223722768536, 1.093305363867116, 1.084388296544119, 1.102065111090516, 1.103875972390869, 1.094857982887685)
    25 = @(1.02, 1.09086299913703, 1.08381326153325, 1.102199918062517, 1.103901202085497, 1.050039783433082, 1.096143429557672, 1.086722169057757, 1.105057419973589, 1.106754041295664, 1.097700078957052)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "{0}{1}" -f $cols[$i], $r
        $ws.Range($cellRef).Value = $vals[$i]
    }
}

Write-Output "Applied 264 cell updates across 24 rows"
